$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'31.003.39"
$ws.Range("E2").Value = "  +0.77%  "

$ws.Range("D3").Value = "'1.958.06"
$ws.Range("E3").Value = "  -0.23%  "

$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").Value = "'244.61"
$ws.Range("E5").Value = "  -1.47%  "

$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").Value = "'0.4868"
$ws.Range("E7").Value = "  +0.65%  "

$ws.Range("D8").Value = "'0.2950"
$ws.Range("E8").Value = "  +0.51%  "

$ws.Range("D9").Value = "'0.06975"
$ws.Range("E9").Value = "  +3.21%  "

$ws.Range("D10").Value = "'19.51"
$ws.Range("E10").Value = "  +2.63%  "

$ws.Range("D11").Value = "'107.75"
$ws.Range("E11").Value = "  -1.23%  "

$ws.Range("D12").Value = "'1.951.60"
$ws.Range("E12").Value = "  -0.46%  "

$ws.Range("D13").Value = "'0.07802"
$ws.Range("E13").Value = "  +0.78%  "

$ws.Range("D14").Value = "'5.489"
$ws.Range("E14").Value = "  +0.23%  "

$ws.Range("D15").Value = "'0.7004"
$ws.Range("E15").Value = "  +0.50%  "

$ws.Range("D16").Value = "'280.85"
$ws.Range("E16").Value = "  -3.82%  "

$ws.Range("D17").Value = "'31.030.33"
$ws.Range("E17").Value = "  +0.88%  "

$ws.Range("D18").Value = "'13.29"
$ws.Range("E18").Value = "  +1.04%  "

$ws.Range("D19").Value = "'0.000007783"
$ws.Range("E19").Value = "  +1.18%  "

$ws.Range("D20").Value = "'2.214.18"
$ws.Range("E20").Value = "  -0.23%  "

$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("D22").Value = "'5.541"
$ws.Range("E22").Value = "  -2.27%  "

$ws.Range("E23").Value = "  -0.38%  "

$ws.Range("D24").Value = "'6.518"
$ws.Range("E24").Value = "  -1.46%  "

$ws.Range("D25").Value = "'9.850"
$ws.Range("E25").Value = "  -0.44%  "

$ws.Range("D26").Value = "'168.42"
$ws.Range("E26").Value = "  -1.06%  "

$ws.Range("D27").Value = "'19.94"
$ws.Range("E27").Value = "  -0.74%  "

$ws.Range("D28").Value = "'2.196"
$ws.Range("E28").Value = "  +0.96%  "

$ws.Range("D29").Value = "'0.1049"
$ws.Range("E29").Value = "  -2.54%  "

$ws.Range("D30").Value = "'1.386"
$ws.Range("E30").Value = "  -4.20%  "

$ws.Range("D31").Value = "'4.630"
$ws.Range("E31").Value = "  -5.15%  "

$ws.Range("E32").Value = "  -2.08%  "

$ws.Range("D33").Value = "'4.463"
$ws.Range("E33").Value = "  -0.45%  "

$ws.Range("D34").Value = "'0.04923"
$ws.Range("E34").Value = "  -3.67%  "

$ws.Range("D35").Value = "'0.7543"
$ws.Range("E35").Value = "  -2.40%  "

$ws.Range("D36").Value = "'1.171"
$ws.Range("E36").Value = "  -0.63%  "

$ws.Range("D37").Value = "'2.733"
$ws.Range("E37").Value = "  +0.32%  "

$ws.Range("D38").Value = "'0.02009"
$ws.Range("E38").Value = "  -1.75%  "

$ws.Range("D39").Value = "'2.694"
$ws.Range("E39").Value = "  -1.17%  "

$ws.Range("D40").Value = "'6.550"
$ws.Range("E40").Value = "  +0.79%  "

$ws.Range("D41").Value = "'77.96"
$ws.Range("E41").Value = "  +11.28%  "

$ws.Range("D42").Value = "'2.141"
$ws.Range("E42").Value = "  +0.52%  "

$ws.Range("D43").Value = "'0.9023"
$ws.Range("E43").Value = "  +1.49%  "

$ws.Range("D44").Value = "'109.24"
$ws.Range("E44").Value = "  -0.70%  "

$ws.Range("D45").Value = "'0.4457"
$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("D46").Value = "'8.122"
$ws.Range("E46").Value = "  +8.57%  "

$ws.Range("D48").Value = "'1.008.60"
$ws.Range("E48").Value = "  +8.05%  "

$ws.Range("D49").Value = "'9.375"
$ws.Range("E49").Value = "  +0.21%  "

$ws.Range("D50").Value = "'0.1255"
$ws.Range("E50").Value = "  -1.66%  "

$ws.Range("D51").Value = "'35.94"
$ws.Range("E51").Value = "  -0.29%  "
